# Daily attendance processing - 2026-01-17 23:58:01
# Swap the order of "System" and the recorder's email address in column G
# ("Recorded By") wherever the cell currently reads "System, <email>".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colG = $ws.Range("G1:G" + $lastRow)
$colG.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
